# smartvac-fhir xlsx mapping updated
#
# - A1 header text: "WHO SmartVacc Parameter" -> "SmartVacc Parameter"
# - E11 note text: "id also present in source – better?"
#                -> "id also present in source – id better?"
# - Selection / scroll position moved from E14 to E12
# - Column D width nudged slightly wider

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text changes ---
$ws.Range("A1").Value = "SmartVacc Parameter"
$ws.Range("E11").Value = "id also present in source – id better?"

# --- Column width tweak (col D, ~63.37 chars) ---
$ws.Columns.Item(4).ColumnWidth = 62.5

# --- View changes: active cell / selection moved to E12 ---
$ws.Range("E12").Select() | Out-Null
